$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: the DeleteUser/Param row - update the value from "Robin" to the
# API user's e-mail address.
$ws.Range("G4").Value = "robin@crate.com"

# Row 5: add a new "Param" row supplying the UserPassword for the API user.
$ws.Range("B5").Value = "Param"
$ws.Range("E5").Value = "UserPassword"
$ws.Range("F5").Value = "string"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "123456"
